$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- Remove the trailing rows (old 24-28) that no longer exist in the new layout ---
$ws.Range("A24:A28").EntireRow.Delete()

# --- Clear cells whose content disappears completely in the new layout (also removes the cell node) ---
$ws.Range("B13:C13").Clear()
$ws.Range("B15:C16").Clear()
$ws.Range("A22:A23").Clear()

# --- Apply correct column formatting (style) to cells that are newly introduced in rows 12-20 ---
# Column A cells (bold, style used by existing single-label rows, e.g. row 11 "Objectives:")
$ws.Cells.Item(11,1).Copy()
$ws.Cells.Item(13,1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(14,1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(15,1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(16,1).PasteSpecial($xlPasteFormats)

# Column B/C cells (style used by existing value rows, e.g. row 14 before edits had B/C with correct style)
$ws.Cells.Item(14,2).Copy()
$ws.Cells.Item(12,2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(18,2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(20,2).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(14,3).Copy()
$ws.Cells.Item(12,3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(18,3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(20,3).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Row 10: Objetivos: / first professor ---
$ws.Cells.Item(10,1).Value = "Objetivos:"
$ws.Cells.Item(10,2).Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Cells.Item(10,3).Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Row 11: Objectives: (unchanged) ---
$ws.Cells.Item(11,1).Value = "Objectives:"

# --- Row 12: Programa resumido: / second professor ---
$ws.Cells.Item(12,1).Value = "Programa resumido:"
$ws.Cells.Item(12,2).Value = "3682251 - Gabrielle Weber Martins"
$ws.Cells.Item(12,3).Value = "3682251 - Gabrielle Weber Martins"

# --- Row 13: Short syllabus: ---
$ws.Cells.Item(13,1).Value = "Short syllabus:"

# --- Row 14: Programa: / third professor ---
$ws.Cells.Item(14,1).Value = "Programa:"
$ws.Cells.Item(14,2).Value = "7797767 - Viktor Pastoukhov"
$ws.Cells.Item(14,3).Value = "7797767 - Viktor Pastoukhov"

# --- Row 15: Syllabus: ---
$ws.Cells.Item(15,1).Value = "Syllabus:"

# --- Row 16: Avaliação: ---
$ws.Cells.Item(16,1).Value = "Avaliação:"

# --- Row 17: Método: / fourth professor ---
$ws.Cells.Item(17,1).Value = "Método:"
$ws.Cells.Item(17,2).Value = "5729033 - Weiliang Qian"
$ws.Cells.Item(17,3).Value = "5729033 - Weiliang Qian"

# --- Row 18: Critério: / evaluation method text ---
$ws.Cells.Item(18,1).Value = "Critério:"
$ws.Cells.Item(18,2).Value = "A avaliação será composta por duas provas escritas (P1 e P2)."
$ws.Cells.Item(18,3).Value = "A avaliação será composta por duas provas escritas (P1 e P2)."

# --- Row 19: Norma de recuperação: / grading criteria text ---
$ws.Cells.Item(19,1).Value = "Norma de recuperação:"
$ws.Cells.Item(19,2).Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Cells.Item(19,3).Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."

# --- Row 20: Bibliografia: / recovery norm text ---
$ws.Cells.Item(20,1).Value = "Bibliografia:"
$ws.Cells.Item(20,2).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Cells.Item(20,3).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# --- Row 21: Requisitos: (unchanged) ---
$ws.Cells.Item(21,1).Value = "Requisitos:"

# --- Row 22: first requirement ---
$ws.Cells.Item(22,2).Value = "LOB1004 -  Cálculo II  (Requisito)" + [char]10
$ws.Cells.Item(22,3).Value = "LOB1004 -  Cálculo II  (Requisito)" + [char]10

# --- Row 23: second requirement ---
$ws.Cells.Item(23,2).Value = "LOB1018 -  Física I  (Requisito)" + [char]10
$ws.Cells.Item(23,3).Value = "LOB1018 -  Física I  (Requisito)" + [char]10

# --- Fix row heights to match the new layout ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 120
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30
